# Fix timeout download files
# Update the "état des virements" worksheet: replace row 2's data (KHADIJA LALA),
# rebuild row 3 as a second CHARIJI ABDELLAH line, and append a new
# ACHENGLI LAILA row plus a totals row, growing the sheet from A1:K3 to A1:K5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: CHARIJI ABDELLAH -------------------------------------------------
$ws.Range("A2").Value = "CHARIJI ABDELLAH"
$ws.Range("B2").Value = "BJ36877"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "00101211111292695000201732"
$ws.Range("D2").Value = "AOURIR"
$ws.Range("E2").Value = "BP CENTRE SUD"
$ws.Range("F2").Value = "Logement de fonction"
$ws.Range("G2").Value = "901/LF/FES "
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 12000
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 10800

# --- Row 3: CHARIJI ABDELLAH (second line, different tax/net) ---------------
$ws.Range("A3").Value = "CHARIJI ABDELLAH"
$ws.Range("B3").Value = "BJ36877"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "00101211111292695000201732"
$ws.Range("D3").Value = "AOURIR"
$ws.Range("E3").Value = "BP CENTRE SUD"
$ws.Range("F3").Value = "Logement de fonction"
$ws.Range("G3").Value = "901/LF/FES "
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 12000
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 11400

# --- Row 4: ACHENGLI LAILA (new row) -----------------------------------------
$ws.Range("A4").Value = "ACHENGLI LAILA"
$ws.Range("B4").Value = "J207703"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "00101211115087750001201090"
$ws.Range("D4").Value = "Ait souss"
$ws.Range("E4").Value = "BP Centre Sud"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "901/FES "
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 10000
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 9500

# --- Row 5: totals row (new row) ---------------------------------------------
$ws.Range("A5").Value = " "
$ws.Range("B5").Value = " "
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = " "
$ws.Range("G5").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("I5").Value = 34000
$ws.Range("J5").Value = 2300
$ws.Range("K5").Value = 31700

# Keep the "number stored as text" error-checking hint in sync with the new
# A1:K5 extent (best effort through the exposed error-checking API).
$ws.Range("A1:K5").Errors.Item(3).Ignore = $true
